$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$co = $ws.ChartObjects().Add(100, 100, 300, 200)
$chart = $co.Chart
$chart.ChartType = -4169
$ser = $chart.SeriesCollection.NewSeries()
$ser.XValues = $ws.Range("BA44:BA64")
$ser.Values = $ws.Range("BB44:BB64")
